# Update the "想去人数" (want-to-go count) column F values on both the
# "展览" and "全部类型" sheets, which carry duplicate data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new value for column F
$updates = @{
    2 = 1338
    3 = 1877
    4 = 167
    5 = 68
    6 = 6311
    7 = 180
    8 = 110
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
